$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Items"
$ws.Range("B1").Value = "Status"
$ws.Range("A2").Value = "testing,shopping,owrjdfnd,43545#@@,party!!"
$ws.Range("B2").Value = "complete,in_progress,complete,complete,in_progress"

$ws.Range("B10").Select()
